$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MD20000.15-DEc")
$ws.Activate()

# New payment entry for row 3: Issue Date (B3), Loan Amount (C3), No of Payments (D3)
# Use the raw date serial number so the existing date number-format/style on
# the cell is preserved instead of Excel re-guessing a new one.
$ws.Range("B3").Value = 45280
$ws.Range("C3").Value = 1600
$ws.Range("D3").Value = 1

# Update the selection to reflect where the user ended up (D3)
$ws.Range("D3").Select()

$wb.Save()
